$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (losing a significant trailing zero), so force them to stay plain text.
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"

$ws.Range("D2").Value = '67.700.41'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '3.783.37'
$ws.Range("E3").Value = '  -1.97%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '597.06'
$ws.Range("E5").Value = '  -0.95%  '
$ws.Range("D6").Value = '169.61'
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("D7").Value = '3.783.65'
$ws.Range("E7").Value = '  -1.95%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.524'
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("E10").Value = '  -0.70%  '
$ws.Range("D11").Value = '6.50'
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("E12").Value = '  -1.08%  '
$ws.Range("D13").Value = '0.0000280'
$ws.Range("E13").Value = '  +4.86%  '
$ws.Range("D14").Value = '36.56'
$ws.Range("E14").Value = '  -1.55%  '
$ws.Range("D15").Value = '4.420.13'
$ws.Range("E15").Value = '  -1.89%  '
$ws.Range("D16").Value = '3.789.60'
$ws.Range("E16").Value = '  -1.87%  '
$ws.Range("D17").Value = '18.55'
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").Value = '67.716.95'
$ws.Range("E18").Value = '  -1.42%  '
$ws.Range("D19").Value = '7.18'
$ws.Range("E19").Value = '  -2.87%  '
$ws.Range("D21").Value = '10.54'
$ws.Range("E21").Value = '  -7.04%  '
$ws.Range("D22").Value = '468.93'
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("E23").Value = '  -2.03%  '
$ws.Range("D24").Value = '0.0000149'
$ws.Range("E24").Value = '  -7.61%  '
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("E26").Value = '  -1.51%  '
$ws.Range("D27").Value = '12.17'
$ws.Range("E27").Value = '  +0.18%  '
$ws.Range("D28").Value = '10.30'
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("E30").Value = '  -1.75%  '
$ws.Range("D31").Value = '3.934.64'
$ws.Range("E31").Value = '  -1.87%  '
$ws.Range("D32").Value = '7.62'
$ws.Range("E32").Value = '  -0.89%  '
$ws.Range("D33").Value = '30.55'
$ws.Range("E33").Value = '  -3.08%  '
$ws.Range("D34").Value = '2.22'
$ws.Range("E34").Value = '  -3.90%  '
$ws.Range("D35").Value = '9.15'
$ws.Range("E35").Value = '  -2.22%  '
$ws.Range("D36").Value = '3.745.99'
$ws.Range("E36").Value = '  -2.03%  '
$ws.Range("D37").Value = '3.80'
$ws.Range("E37").Value = '  +2.29%  '
$ws.Range("E38").Value = '  -0.53%  '
$ws.Range("E39").Value = '  -1.20%  '
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("E41").Value = '  -2.37%  '
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("E43").Value = '  -1.09%  '
$ws.Range("E45").Value = '  -0.56%  '
$ws.Range("E46").Value = '  -2.03%  '
$ws.Range("D47").Value = '45.81'
$ws.Range("E47").Value = '  -2.69%  '
$ws.Range("D48").Value = '396.05'
$ws.Range("E48").Value = '  -5.08%  '
$ws.Range("E49").Value = '  -7.93%  '
$ws.Range("D50").Value = '140.59'
$ws.Range("E50").Value = '  -0.98%  '
$ws.Range("D51").Value = '39.29'
$ws.Range("E51").Value = '  +3.36%  '
